# BUG Pseudorandomized Con Schedules
# Making sure the Neg without UCS is never shown first
#
# 1) Remove the bogus empty placeholder cells (D:G, rows 2-18) on the
#    PreCond1 sheet - these were artifacts of the earlier (buggy)
#    pseudo-randomization run and carried no real data.
# 2) Add two new sheets - PreCond2 and RatingPreCond1 - containing the
#    results of a re-run with a fixed (non-Neg-first) schedule.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Clean up PreCond1 (sheet7): remove empty D:G placeholders for rows
#    2-18 (row 19 has real trigger data and must stay untouched).
# ---------------------------------------------------------------------
$preCond1 = $wb.Worksheets.Item("PreCond1")
$preCond1.Range("D2:G18").ClearContents()

# ---------------------------------------------------------------------
# 2) Add "PreCond2" sheet after the last existing sheet (PreCond1)
# ---------------------------------------------------------------------
$preCond2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$preCond2.Name = "PreCond2"

$preCond2.Range("A1").Value = "PreCondName"
$preCond2.Range("B1").Value = "TrgCol"
$preCond2.Range("C1").Value = "n"
$preCond2.Range("D1").Value = "PreCondTriggKey.keys_raw"
$preCond2.Range("E1").Value = "PreCondTriggKey.rt_mean"
$preCond2.Range("F1").Value = "PreCondTriggKey.rt_raw"
$preCond2.Range("G1").Value = "PreCondTriggKey.rt_std"
$preCond2.Range("H1").Value = "PreCondTriggMouse.leftButton_raw"
$preCond2.Range("I1").Value = "PreCondTriggMouse.midButton_raw"
$preCond2.Range("J1").Value = "PreCondTriggMouse.rightButton_raw"
$preCond2.Range("K1").Value = "PreCondTriggMouse.time_raw"
$preCond2.Range("L1").Value = "PreCondTriggMouse.x_raw"
$preCond2.Range("M1").Value = "PreCondTriggMouse.y_raw"
$preCond2.Range("N1").Value = "order"

$preCond2Data = @(
    @(2,  "Stimuli/Neg.BMP",  0,   8),
    @(3,  "Stimuli/Neg.BMP",  0,   16),
    @(4,  "Stimuli/Neg.BMP",  0,   6),
    @(5,  "Stimuli/Neg.BMP",  0,   10),
    @(6,  "Stimuli/Neg.BMP",  0,   12),
    @(7,  "Stimuli/Neu.BMP",  0.4, 5),
    @(8,  "Stimuli/Neu.BMP",  0.4, 0),
    @(9,  "Stimuli/Neu.BMP",  0.4, 1),
    @(10, "Stimuli/Neu.BMP",  0.4, 13),
    @(11, "Stimuli/Neu.BMP",  0.4, 4),
    @(12, "Stimuli/Pos.BMP",  0.6, 14),
    @(13, "Stimuli/Pos.BMP",  0.6, 17),
    @(14, "Stimuli/Pos.BMP",  0.6, 3),
    @(15, "Stimuli/Pos.BMP",  0.6, 15),
    @(16, "Stimuli/Pos.BMP",  0.6, 9),
    @(18, "Stimuli/Trig.BMP", 0.9, 2)
)

foreach ($row in $preCond2Data) {
    $r = $row[0]
    $preCond2.Range("A$r").Value = $row[1]
    $preCond2.Range("B$r").Value = $row[2]
    $preCond2.Range("C$r").Value = 1
    $preCond2.Range("H$r").Value = $row[3]
}

# Row 17 and 19 - Trig rows that had an actual key-press response
$preCond2.Range("A17").Value = "Stimuli/Trig.BMP"
$preCond2.Range("B17").Value = 0.9
$preCond2.Range("C17").Value = 1
$preCond2.Range("D17").Value = "''space'"
$preCond2.Range("E17").Value = 1.07956600189209
$preCond2.Range("F17").Value = 1.07956600189209
$preCond2.Range("G17").Value = 0
$preCond2.Range("H17").Value = 11

$preCond2.Range("A19").Value = "Stimuli/Trig.BMP"
$preCond2.Range("B19").Value = 0.9
$preCond2.Range("C19").Value = 1
$preCond2.Range("D19").Value = "''space'"
$preCond2.Range("E19").Value = 5.133314609527588
$preCond2.Range("F19").Value = 5.133314609527588
$preCond2.Range("G19").Value = 0
$preCond2.Range("H19").Value = 7

$preCond2.Range("A21").Value = "extraInfo"

$preCond2.Range("A22").Value = "Participant_ID"
$preCond2.Range("B22").Value = "S00"
$preCond2.Range("A23").Value = "Version"
$preCond2.Range("B23").Value = "Short"
$preCond2.Range("A24").Value = "Language"
$preCond2.Range("B24").Value = "EN"
$preCond2.Range("A25").Value = "date"
$preCond2.Range("B25").Value = "2023-05-28_15h36.27.014"
$preCond2.Range("A26").Value = "expName"
$preCond2.Range("B26").Value = "TCET"
$preCond2.Range("A27").Value = "psychopyVersion"
$preCond2.Range("B27").NumberFormat = "@"
$preCond2.Range("B27").Value = "2023.1.1"
$preCond2.Range("A28").Value = "frameRate"
$preCond2.Range("B28").Value = 60.40806857672084

# ---------------------------------------------------------------------
# 3) Add "RatingPreCond1" sheet after PreCond2
# ---------------------------------------------------------------------
$ratingPreCond1 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ratingPreCond1.Name = "RatingPreCond1"

$ratingPreCond1.Range("A1").Value = "CondName"
$ratingPreCond1.Range("B1").Value = "TrgCol"
$ratingPreCond1.Range("C1").Value = "FearRatingInstTextEN"
$ratingPreCond1.Range("D1").Value = "FearRatingInstTextCN"
$ratingPreCond1.Range("E1").Value = "n"
$ratingPreCond1.Range("F1").Value = "RatingKey.keys_raw"
$ratingPreCond1.Range("G1").Value = "RatingKey.rt_mean"
$ratingPreCond1.Range("H1").Value = "RatingKey.rt_raw"
$ratingPreCond1.Range("I1").Value = "RatingKey.rt_std"
$ratingPreCond1.Range("J1").Value = "Scale.response_mean"
$ratingPreCond1.Range("K1").Value = "Scale.response_raw"
$ratingPreCond1.Range("L1").Value = "Scale.response_std"
$ratingPreCond1.Range("M1").Value = "Scale.rt_mean"
$ratingPreCond1.Range("N1").Value = "Scale.rt_raw"
$ratingPreCond1.Range("O1").Value = "Scale.rt_std"
$ratingPreCond1.Range("P1").Value = "Submit.numClicks_mean"
$ratingPreCond1.Range("Q1").Value = "Submit.numClicks_raw"
$ratingPreCond1.Range("R1").Value = "Submit.numClicks_std"
$ratingPreCond1.Range("S1").Value = "Submit.timesOff_raw"
$ratingPreCond1.Range("T1").Value = "Submit.timesOn_raw"
$ratingPreCond1.Range("U1").Value = "mouse.leftButton_raw"
$ratingPreCond1.Range("V1").Value = "mouse.midButton_raw"
$ratingPreCond1.Range("W1").Value = "mouse.rightButton_raw"
$ratingPreCond1.Range("X1").Value = "mouse.time_raw"
$ratingPreCond1.Range("Y1").Value = "mouse.x_raw"
$ratingPreCond1.Range("Z1").Value = "mouse.y_raw"
$ratingPreCond1.Range("AA1").Value = "order"

# Row 2 - Stimuli/Raw_Neg.BMP
$ratingPreCond1.Range("A2").Value = "Stimuli/Raw_Neg.BMP"
$ratingPreCond1.Range("B2").Value = 0
$ratingPreCond1.Range("C2").Value = "How do you feel about this face?"
$ratingPreCond1.Range("D2").Value = "你对这个面孔有什么感觉?（'– '表示负面恐怖的, '+'表示积极快乐的）"
$ratingPreCond1.Range("E2").Value = 1
$ratingPreCond1.Range("F2").Value = "''space'"
$ratingPreCond1.Range("G2").Value = 0.3672721982002258
$ratingPreCond1.Range("H2").Value = 0.3672721982002258
$ratingPreCond1.Range("I2").Value = 0
$ratingPreCond1.Range("J2").Value = 10
$ratingPreCond1.Range("K2").Value = 10
$ratingPreCond1.Range("L2").Value = 0
$ratingPreCond1.Range("M2").Value = 1.187361001968384
$ratingPreCond1.Range("N2").Value = 1.187361001968384
$ratingPreCond1.Range("O2").Value = 0
$ratingPreCond1.Range("P2").Value = 0
$ratingPreCond1.Range("Q2").Value = 0
$ratingPreCond1.Range("R2").Value = 0
$ratingPreCond1.Range("S2").Value = "'''"
$ratingPreCond1.Range("T2").Value = "'''"
$ratingPreCond1.Range("U2").Value = 1
$ratingPreCond1.Range("V2").Value = 0
$ratingPreCond1.Range("W2").Value = 0
$ratingPreCond1.Range("X2").Value = 1.062335400027223
$ratingPreCond1.Range("Y2").Value = 0.6030092592592593
$ratingPreCond1.Range("Z2").Value = -0.5120370370370371
$ratingPreCond1.Range("AA2").Value = 0

# Row 3 - Stimuli/Raw_Neu.BMP
$ratingPreCond1.Range("A3").Value = "Stimuli/Raw_Neu.BMP"
$ratingPreCond1.Range("B3").Value = 0.4
$ratingPreCond1.Range("E3").Value = 1
$ratingPreCond1.Range("F3").Value = "''space'"
$ratingPreCond1.Range("G3").Value = 0.1364340931177139
$ratingPreCond1.Range("H3").Value = 0.1364340931177139
$ratingPreCond1.Range("I3").Value = 0
$ratingPreCond1.Range("J3").Value = 10
$ratingPreCond1.Range("K3").Value = 10
$ratingPreCond1.Range("L3").Value = 0
$ratingPreCond1.Range("M3").Value = 0.3078139126300812
$ratingPreCond1.Range("N3").Value = 0.3078139126300812
$ratingPreCond1.Range("O3").Value = 0
$ratingPreCond1.Range("P3").Value = 0
$ratingPreCond1.Range("Q3").Value = 0
$ratingPreCond1.Range("R3").Value = 0
$ratingPreCond1.Range("S3").Value = "'''"
$ratingPreCond1.Range("T3").Value = "'''"
$ratingPreCond1.Range("U3").Value = 1
$ratingPreCond1.Range("V3").Value = 0
$ratingPreCond1.Range("W3").Value = 0
$ratingPreCond1.Range("X3").Value = 0.1759938999894075
$ratingPreCond1.Range("Y3").Value = 0.6030092592592593
$ratingPreCond1.Range("Z3").Value = -0.5120370370370371
$ratingPreCond1.Range("AA3").Value = 1

# Row 4 - Stimuli/Raw_Pos.BMP
$ratingPreCond1.Range("A4").Value = "Stimuli/Raw_Pos.BMP"
$ratingPreCond1.Range("B4").Value = 0.6
$ratingPreCond1.Range("E4").Value = 1
$ratingPreCond1.Range("F4").Value = "''space'"
$ratingPreCond1.Range("G4").Value = 0.08793230354785919
$ratingPreCond1.Range("H4").Value = 0.08793230354785919
$ratingPreCond1.Range("I4").Value = 0
$ratingPreCond1.Range("J4").Value = 10
$ratingPreCond1.Range("K4").Value = 10
$ratingPreCond1.Range("L4").Value = 0
$ratingPreCond1.Range("M4").Value = 0.2875483930110931
$ratingPreCond1.Range("N4").Value = 0.2875483930110931
$ratingPreCond1.Range("O4").Value = 0
$ratingPreCond1.Range("P4").Value = 0
$ratingPreCond1.Range("Q4").Value = 0
$ratingPreCond1.Range("R4").Value = 0
$ratingPreCond1.Range("S4").Value = "'''"
$ratingPreCond1.Range("T4").Value = "'''"
$ratingPreCond1.Range("U4").Value = 1
$ratingPreCond1.Range("V4").Value = 0
$ratingPreCond1.Range("W4").Value = 0
$ratingPreCond1.Range("X4").Value = 0.2114670000155456
$ratingPreCond1.Range("Y4").Value = 0.6030092592592593
$ratingPreCond1.Range("Z4").Value = -0.5120370370370371
$ratingPreCond1.Range("AA4").Value = 2

# Row 5 - Stimuli/Raw_Trig.BMP
$ratingPreCond1.Range("A5").Value = "Stimuli/Raw_Trig.BMP"
$ratingPreCond1.Range("B5").Value = 0.9
$ratingPreCond1.Range("E5").Value = 1
$ratingPreCond1.Range("F5").Value = "''space'"
$ratingPreCond1.Range("G5").Value = 0.3338384032249451
$ratingPreCond1.Range("H5").Value = 0.3338384032249451
$ratingPreCond1.Range("I5").Value = 0
$ratingPreCond1.Range("J5").Value = 10
$ratingPreCond1.Range("K5").Value = 10
$ratingPreCond1.Range("L5").Value = 0
$ratingPreCond1.Range("M5").Value = 0.472949206829071
$ratingPreCond1.Range("N5").Value = 0.472949206829071
$ratingPreCond1.Range("O5").Value = 0
$ratingPreCond1.Range("P5").Value = 0
$ratingPreCond1.Range("Q5").Value = 0
$ratingPreCond1.Range("R5").Value = 0
$ratingPreCond1.Range("S5").Value = "'''"
$ratingPreCond1.Range("T5").Value = "'''"
$ratingPreCond1.Range("U5").Value = 1
$ratingPreCond1.Range("V5").Value = 0
$ratingPreCond1.Range("W5").Value = 0
$ratingPreCond1.Range("X5").Value = 0.4037122000008821
$ratingPreCond1.Range("Y5").Value = 0.6030092592592593
$ratingPreCond1.Range("Z5").Value = -0.5120370370370371
$ratingPreCond1.Range("AA5").Value = 3

$ratingPreCond1.Range("A7").Value = "extraInfo"

$ratingPreCond1.Range("A8").Value = "Participant_ID"
$ratingPreCond1.Range("B8").Value = "S00"
$ratingPreCond1.Range("A9").Value = "Version"
$ratingPreCond1.Range("B9").Value = "Short"
$ratingPreCond1.Range("A10").Value = "Language"
$ratingPreCond1.Range("B10").Value = "EN"
$ratingPreCond1.Range("A11").Value = "date"
$ratingPreCond1.Range("B11").Value = "2023-05-28_15h36.27.014"
$ratingPreCond1.Range("A12").Value = "expName"
$ratingPreCond1.Range("B12").Value = "TCET"
$ratingPreCond1.Range("A13").Value = "psychopyVersion"
$ratingPreCond1.Range("B13").NumberFormat = "@"
$ratingPreCond1.Range("B13").Value = "2023.1.1"
$ratingPreCond1.Range("A14").Value = "frameRate"
$ratingPreCond1.Range("B14").Value = 60.40806857672084
